$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking values
# (e.g. "1.00", "309.28") are written as literal text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "40.045.36"
$ws.Cells.Item(2, 5).Value = "  -2.98%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.337.53"
$ws.Cells.Item(3, 5).Value = "  -4.07%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  -0.16%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "309.28"
$ws.Cells.Item(5, 5).Value = "  -2.33%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "85.16"
$ws.Cells.Item(6, 5).Value = "  -4.85%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -2.30%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.01%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -2.38%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -1.62%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "30.15"
$ws.Cells.Item(11, 5).Value = "  -5.99%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +1.17%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "2.696.83"
$ws.Cells.Item(13, 5).Value = "  -4.19%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -3.83%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "14.73"
$ws.Cells.Item(15, 5).Value = "  -3.71%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.357.98"
$ws.Cells.Item(16, 5).Value = "  -3.32%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "0.758"
$ws.Cells.Item(17, 5).Value = "  -1.56%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "40.018.95"
$ws.Cells.Item(18, 5).Value = "  -2.95%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "0.0₃0902"
$ws.Cells.Item(19, 5).Value = "  -1.83%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -1.77%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "67.98"
$ws.Cells.Item(21, 5).Value = "  -5.18%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "10.67"
$ws.Cells.Item(22, 5).Value = "  -3.36%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "235.73"
$ws.Cells.Item(23, 5).Value = "  +0.64%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "2.56"
$ws.Cells.Item(24, 5).Value = "  -5.04%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.11%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "1.82"
$ws.Cells.Item(26, 5).Value = "  -3.02%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -2.61%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -4.26%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "9.29"
$ws.Cells.Item(29, 5).Value = "  -2.47%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "34.57"
$ws.Cells.Item(30, 5).Value = "  -0.17%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "153.81"
$ws.Cells.Item(31, 5).Value = "  -2.20%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -0.11%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "5.11"
$ws.Cells.Item(33, 5).Value = "  -2.63%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -3.77%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -3.33%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.40%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "2.76"
$ws.Cells.Item(37, 5).Value = "  -4.70%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.0992"
$ws.Cells.Item(38, 5).Value = "  -0.11%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Celestia"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(39, 4).Value = "15.61"
$ws.Cells.Item(39, 5).Value = "  -5.50%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "ARBITRUM"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(40, 4).Value = "1.73"
$ws.Cells.Item(40, 5).Value = "  -2.32%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "3.86"
$ws.Cells.Item(41, 5).Value = "  -0.49%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.950.08"
$ws.Cells.Item(42, 5).Value = "  -1.65%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -4.80%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "17.73"
$ws.Cells.Item(44, 5).Value = "  -1.70%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -4.23%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "9.44"
$ws.Cells.Item(46, 5).Value = "  +0.12%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -4.91%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "2.554.74"
$ws.Cells.Item(48, 5).Value = "  -4.51%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "92.88"
$ws.Cells.Item(49, 5).Value = "  -2.30%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "70.71"
$ws.Cells.Item(50, 5).Value = "  -3.10%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "50.09"
$ws.Cells.Item(51, 5).Value = "  -3.25%  "

# Restore column D to its original (General) formatting now that the
# text values have been written, so no extra number format lingers.
$ws.Range("D2:D51").ClearFormats()
